# Generate Report for Handoff
# Rename the localized markdown/xliff artifacts from the old GUID
# (e02418ef-5a14-4557-ac77-68694f6e3506) to the new one
# (c9a46d89-45ba-4db9-838f-7659a4255cb9) across all three sheets, and
# bump the related handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "e02418ef-5a14-4557-ac77-68694f6e3506"
$newGuid = "c9a46d89-45ba-4db9-838f-7659a4255cb9"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7bdd2368fee7c098b67c8fc4d884ae836d3f1e6/e2e/"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
# A2: File Name
$ws1.Range("A2").Value = "$newGuid.md"

# B2: Path And Name (hyperlink display text + target)
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "$ghBase$newGuid.md", "", "", "e2e\$newGuid.md")

# G2: Latest HO Xliff Generate Date
$ws1.Range("G2").Value = "2016-08-19 02:54:34"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
# A2: Source File Name (hyperlink display text + target)
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "$ghBase$newGuid.md", "", "", "$newGuid.md")

# G2: Latest Handoff File
$ws2.Range("G2").Value = "$newGuid.94665bda437ee677dc4f3e8b9f53d435807ebe71.zh-cn.xlf"

# H2: Latest Handoff Datetime
$ws2.Range("H2").Value = "2016-08-19 02:54:29"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
# A2: Source File Name (hyperlink display text + target)
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "$ghBase$newGuid.md", "", "", "$newGuid.md")

# G2: Latest Handoff File
$ws3.Range("G2").Value = "$newGuid.94665bda437ee677dc4f3e8b9f53d435807ebe71.de-de.xlf"

# H2: Latest Handoff Datetime
# (shares the same underlying text as the Overview sheet's "Latest HO Xliff
# Generate Date" in the original workbook, so it moves in lock-step with it)
$ws3.Range("H2").Value = "2016-08-19 02:54:34"
